$p = $ppt.ActivePresentation
$sm = $p.SlideMaster
$nm = $p.NotesMaster
$sm.Name = "SMNameTest"
Write-Host "NotesMaster name after setting SlideMaster name: $($nm.Name)"
